$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsRM = $wb.Worksheets.Item("RM")

# Remove the "U.S. EPS 3.0.0 Note" block (rows 14-17) from the About sheet
$wsAbout.Range("A14:B17").Clear()

# Update RM sheet row 2 values from 0 to 0.1412 (B2:AK2)
$wsRM.Range("B2:AK2").Value = 0.1412
